# "R script: all grouping done by unique_identifier"
#
# - end_message (H2) gets an extra instruction line.
# - accuracyCriterion (I2/I3) raised from 0 to 80.
# - n_pairs_test_blocks (R2/R3) raised from 1 to 3.
# - column Y widened to fit the longer end_message.
# - the view's selection moves to H3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. end_message (H2): append a second line instructing the user to press Enter.
$ws.Range("H2").Value = "End of task 1. Take a short break.`nPress the enter key to continue."

# 2. accuracyCriterion (column I), rows 2 and 3: 0 -> 80
$ws.Range("I2").Value = 80
$ws.Range("I3").Value = 80

# 3. n_pairs_test_blocks (column R), rows 2 and 3: 1 -> 3
$ws.Range("R2").Value = 3
$ws.Range("R3").Value = 3

# 4. Widen column Y (25) so the longer end_message text is readable.
$ws.Columns.Item(25).ColumnWidth = 20.428571428571427

# 5. Move the selection to H3 (matches the recorded sheet view).
$ws.Range("H3").Select()
